$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Simple Sheet")
$ws.Rows.Item(1).RowHeight = 20.25
$ws.Rows.Item(2).RowHeight = 15.75
